# Applies odds updates for Jogos_da_Semana_FlashScore_2024-10-22.xlsx
# as described in the commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Cells.Item(7, 7).Value = 1.85   # G7: 1.8 -> 1.85
$ws.Cells.Item(7, 9).Value = 5   # I7: 5.25 -> 5
$ws.Cells.Item(7, 13).Value = 1.11   # M7: 1.13 -> 1.11
$ws.Cells.Item(7, 14).Value = 6.5   # N7: 6 -> 6.5
$ws.Cells.Item(7, 19).Value = 1.62   # S7: 1.57 -> 1.62
$ws.Cells.Item(7, 20).Value = 2.2   # T7: 2.25 -> 2.2
$ws.Cells.Item(7, 27).Value = 21   # AA7: 19 -> 21
$ws.Cells.Item(7, 32).Value = 101   # AF7: 81 -> 101
$ws.Cells.Item(7, 46).Value = 2.2   # AT7: 2.25 -> 2.2

# Row 8
$ws.Cells.Item(8, 17).Value = 2.4   # Q8: 2.35 -> 2.4
$ws.Cells.Item(8, 18).Value = 1.53   # R8: 1.57 -> 1.53

# Row 37
$ws.Cells.Item(37, 7).Value = 4.1   # G37: 5 -> 4.1
$ws.Cells.Item(37, 8).Value = 3.6   # H37: 3.75 -> 3.6
$ws.Cells.Item(37, 9).Value = 1.85   # I37: 1.7 -> 1.85
$ws.Cells.Item(37, 10).Value = 4.33   # J37: 5 -> 4.33
$ws.Cells.Item(37, 11).Value = 2.25   # K37: 2.3 -> 2.25
$ws.Cells.Item(37, 12).Value = 2.5   # L37: 2.25 -> 2.5
$ws.Cells.Item(37, 14).Value = 12   # N37: 13 -> 12
$ws.Cells.Item(37, 17).Value = 1.8   # Q37: 1.73 -> 1.8
$ws.Cells.Item(37, 18).Value = 2   # R37: 2.08 -> 2
$ws.Cells.Item(37, 19).Value = 1.36   # S37: 1.33 -> 1.36
$ws.Cells.Item(37, 20).Value = 3   # T37: 3.25 -> 3
$ws.Cells.Item(37, 21).Value = 1.7   # U37: 1.75 -> 1.7
$ws.Cells.Item(37, 22).Value = 2.05   # V37: 2 -> 2.05
$ws.Cells.Item(37, 23).Value = 13   # W37: 15 -> 13
$ws.Cells.Item(37, 24).Value = 21   # X37: 26 -> 21
$ws.Cells.Item(37, 25).Value = 13   # Y37: 15 -> 13
$ws.Cells.Item(37, 26).Value = 41   # Z37: 51 -> 41
$ws.Cells.Item(37, 27).Value = 29   # AA37: 41 -> 29
$ws.Cells.Item(37, 28).Value = 34   # AB37: 41 -> 34
$ws.Cells.Item(37, 30).Value = 6.5   # AD37: 7 -> 6.5
$ws.Cells.Item(37, 31).Value = 13   # AE37: 15 -> 13
$ws.Cells.Item(37, 33).Value = 151   # AG37: 201 -> 151
$ws.Cells.Item(37, 34).Value = 8   # AH37: 7.5 -> 8
$ws.Cells.Item(37, 35).Value = 9.5   # AI37: 8.5 -> 9.5
$ws.Cells.Item(37, 37).Value = 15   # AK37: 13 -> 15
$ws.Cells.Item(37, 38).Value = 15   # AL37: 13 -> 15
$ws.Cells.Item(37, 40).Value = 6   # AN37: 6.5 -> 6
$ws.Cells.Item(37, 41).Value = 21   # AO37: 23 -> 21
$ws.Cells.Item(37, 42).Value = 26   # AP37: 29 -> 26
$ws.Cells.Item(37, 43).Value = 67   # AQ37: 81 -> 67
$ws.Cells.Item(37, 44).Value = 81   # AR37: 101 -> 81
$ws.Cells.Item(37, 45).Value = 151   # AS37: 201 -> 151
$ws.Cells.Item(37, 46).Value = 3   # AT37: 3.25 -> 3
$ws.Cells.Item(37, 47).Value = 7.5   # AU37: 8 -> 7.5
$ws.Cells.Item(37, 49).Value = 4   # AW37: 3.75 -> 4
$ws.Cells.Item(37, 50).Value = 10   # AX37: 8.5 -> 10
$ws.Cells.Item(37, 52).Value = 34   # AZ37: 29 -> 34
$ws.Cells.Item(37, 53).Value = 51   # BA37: 41 -> 51
$ws.Cells.Item(37, 56).Value = 126   # BD37: 151 -> 126

# Row 38
$ws.Cells.Item(38, 7).Value = 1.75   # G38: 1.7 -> 1.75
$ws.Cells.Item(38, 9).Value = 5   # I38: 5.25 -> 5
$ws.Cells.Item(38, 10).Value = 2.38   # J38: 2.3 -> 2.38
$ws.Cells.Item(38, 15).Value = 1.29   # O38: 1.25 -> 1.29
$ws.Cells.Item(38, 16).Value = 3.5   # P38: 3.75 -> 3.5
$ws.Cells.Item(38, 17).Value = 1.95   # Q38: 1.88 -> 1.95
$ws.Cells.Item(38, 18).Value = 1.9   # R38: 1.98 -> 1.9
$ws.Cells.Item(38, 23).Value = 7.5   # W38: 7 -> 7.5
$ws.Cells.Item(38, 24).Value = 8.5   # X38: 8 -> 8.5
$ws.Cells.Item(38, 29).Value = 10   # AC38: 11 -> 10
$ws.Cells.Item(38, 30).Value = 6.5   # AD38: 7 -> 6.5
$ws.Cells.Item(38, 45).Value = 126   # AS38: 151 -> 126
$ws.Cells.Item(38, 50).Value = 23   # AX38: 26 -> 23
$ws.Cells.Item(38, 51).Value = 29   # AY38: 34 -> 29

# Row 39
$ws.Cells.Item(39, 7).Value = 1.87   # G39: 1.78 -> 1.87
$ws.Cells.Item(39, 9).Value = 3.8   # I39: 4.25 -> 3.8
$ws.Cells.Item(39, 10).Value = 2.5   # J39: 2.37 -> 2.5
$ws.Cells.Item(39, 12).Value = 4.25   # L39: 4.6 -> 4.25
$ws.Cells.Item(39, 13).Value = 1.05   # M39: 8.300000000000001 -> 1.05
$ws.Cells.Item(39, 14).Value = 9   # N39: 1.03 -> 9
$ws.Cells.Item(39, 16).Value = 2.95   # P39: 3 -> 2.95
$ws.Cells.Item(39, 17).Value = 1.93   # Q39: 1.9 -> 1.93
$ws.Cells.Item(39, 18).Value = 1.78   # R39: 1.8 -> 1.78
$ws.Cells.Item(39, 23).Value = 7   # W39: 6.8 -> 7
$ws.Cells.Item(39, 24).Value = 8.75   # X39: 8.25 -> 8.75
$ws.Cells.Item(39, 25).Value = 8.25   # Y39: 8 -> 8.25
$ws.Cells.Item(39, 26).Value = 16   # Z39: 14.5 -> 16
$ws.Cells.Item(39, 27).Value = 15   # AA39: 14 -> 15
$ws.Cells.Item(39, 28).Value = 27   # AB39: 26 -> 27
$ws.Cells.Item(39, 29).Value = 9.75   # AC39: 9.5 -> 9.75
$ws.Cells.Item(39, 30).Value = 6.6   # AD39: 6.7 -> 6.6
$ws.Cells.Item(39, 34).Value = 10.75   # AH39: 11.75 -> 10.75
$ws.Cells.Item(39, 35).Value = 21   # AI39: 24 -> 21
$ws.Cells.Item(39, 36).Value = 13   # AJ39: 14 -> 13
$ws.Cells.Item(39, 37).Value = 55   # AK39: 70 -> 55
$ws.Cells.Item(39, 38).Value = 35   # AL39: 40 -> 35
$ws.Cells.Item(39, 39).Value = 40   # AM39: 50 -> 40
$ws.Cells.Item(39, 40).Value = 3.7   # AN39: 3.6 -> 3.7
$ws.Cells.Item(39, 41).Value = 9.5   # AO39: 9 -> 9.5
$ws.Cells.Item(39, 42).Value = 19   # AP39: 18.5 -> 19
$ws.Cells.Item(39, 43).Value = 35   # AQ39: 32 -> 35
$ws.Cells.Item(39, 44).Value = 70   # AR39: 65 -> 70
$ws.Cells.Item(39, 47).Value = 7.4   # AU39: 7.3 -> 7.4
$ws.Cells.Item(39, 49).Value = 5.5   # AW39: 5.9 -> 5.5
$ws.Cells.Item(39, 50).Value = 21   # AX39: 25 -> 21
$ws.Cells.Item(39, 51).Value = 29   # AY39: 30 -> 29
$ws.Cells.Item(39, 52).Value = 120   # AZ39: 150 -> 120
$ws.Cells.Item(39, 53).Value = 150   # BA39: 175 -> 150

# Row 51
$ws.Cells.Item(51, 17).Value = 1.85   # Q51: 1.88 -> 1.85
$ws.Cells.Item(51, 18).Value = 2   # R51: 1.98 -> 2
